$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$newText = "33 x 98" + [char]11 + "  9    8" + [char]11 + "  ----" + [char]11 + "3|    |" + [char]11 + "3|    |"
$t.Cell(1, 1).Range.Text = $newText

$newText = "81 x 95" + [char]11 + "  9    5" + [char]11 + "  ----" + [char]11 + "8|    |" + [char]11 + "1|    |"
$t.Cell(1, 2).Range.Text = $newText

$newText = "19 x 96" + [char]11 + "  9    6" + [char]11 + "  ----" + [char]11 + "1|    |" + [char]11 + "9|    |"
$t.Cell(1, 3).Range.Text = $newText

$newText = "47 x 62" + [char]11 + "  6    2" + [char]11 + "  ----" + [char]11 + "4|    |" + [char]11 + "7|    |"
$t.Cell(2, 1).Range.Text = $newText

$newText = "58 x 86" + [char]11 + "  8    6" + [char]11 + "  ----" + [char]11 + "5|    |" + [char]11 + "8|    |"
$t.Cell(2, 2).Range.Text = $newText

$newText = "89 x 37" + [char]11 + "  3    7" + [char]11 + "  ----" + [char]11 + "8|    |" + [char]11 + "9|    |"
$t.Cell(2, 3).Range.Text = $newText

$newText = "11 x 62" + [char]11 + "  6    2" + [char]11 + "  ----" + [char]11 + "1|    |" + [char]11 + "1|    |"
$t.Cell(3, 1).Range.Text = $newText

$newText = "71 x 47" + [char]11 + "  4    7" + [char]11 + "  ----" + [char]11 + "7|    |" + [char]11 + "1|    |"
$t.Cell(3, 2).Range.Text = $newText

$newText = "62 x 69" + [char]11 + "  6    9" + [char]11 + "  ----" + [char]11 + "6|    |" + [char]11 + "2|    |"
$t.Cell(3, 3).Range.Text = $newText

$newText = "36 x 92" + [char]11 + "  9    2" + [char]11 + "  ----" + [char]11 + "3|    |" + [char]11 + "6|    |"
$t.Cell(4, 1).Range.Text = $newText

$newText = "61 x 13" + [char]11 + "  1    3" + [char]11 + "  ----" + [char]11 + "6|    |" + [char]11 + "1|    |"
$t.Cell(4, 2).Range.Text = $newText

$newText = "24 x 19" + [char]11 + "  1    9" + [char]11 + "  ----" + [char]11 + "2|    |" + [char]11 + "4|    |"
$t.Cell(4, 3).Range.Text = $newText

$newText = "77 x 79" + [char]11 + "  7    9" + [char]11 + "  ----" + [char]11 + "7|    |" + [char]11 + "7|    |"
$t.Cell(5, 1).Range.Text = $newText

$newText = "97 x 45" + [char]11 + "  4    5" + [char]11 + "  ----" + [char]11 + "9|    |" + [char]11 + "7|    |"
$t.Cell(5, 2).Range.Text = $newText

$newText = "76 x 71" + [char]11 + "  7    1" + [char]11 + "  ----" + [char]11 + "7|    |" + [char]11 + "6|    |"
$t.Cell(5, 3).Range.Text = $newText
